$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Commit: "Add artificial endpoint to titration curves"
#
# For each of the three TOC-class titration blocks (TOC <= 3,
# 3 < TOC <= 5, TOC > 5) append one artificial endpoint row that models a
# saturated solution of CaCO3 at room temperature/pressure: 13 mg/l CaCO3
# at pH 8.4. This constrains downstream model extrapolation to physically
# plausible pH values.
#
# Insert the bottom-most group first so the row numbers of the groups
# still to be processed are not disturbed by the earlier inserts.

# --- Group "TOC > 5" currently occupies rows 44:64 -> append after row 64
$ws.Rows.Item(65).Insert()
$ws.Range("A65").Value2 = "TOC > 5"
$ws.Range("B65").Value2 = 13
$ws.Range("C65").Value2 = 8.4

# --- Group "3 < TOC <= 5" currently occupies rows 23:43 -> append after row 43
$ws.Rows.Item(44).Insert()
$ws.Range("A44").Value2 = "3 < TOC ≤ 5"
$ws.Range("B44").Value2 = 13
$ws.Range("C44").Value2 = 8.4

# --- Group "TOC <= 3" currently occupies rows 2:22 -> append after row 22
$ws.Rows.Item(23).Insert()
$ws.Range("A23").Value2 = "TOC ≤ 3"
$ws.Range("B23").Value2 = 13
$ws.Range("C23").Value2 = 8.4

# Reflect the author's final selection/scroll position (row 32 area,
# cell C68 selected) as closely as the object model allows.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
[void]$ws.Range("A2").Select()
$win.FreezePanes = $true
[void]$ws.Range("C68").Select()
